# Updated cryptos list on Wed Jun  7 13:26:33 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '26.776.94'
$ws.Cells.Item(2, 5).Value = '  +4.83%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.864.78'
$ws.Cells.Item(3, 5).Value = '  +3.13%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.09%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''271.00'
$ws.Cells.Item(5, 5).Value = '  -1.36%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -0.09%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''0.5313'
$ws.Cells.Item(7, 5).Value = '  +6.23%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.3361'
$ws.Cells.Item(8, 5).Value = '  -1.57%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.06820'
$ws.Cells.Item(9, 5).Value = '  +3.01%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''19.80'
$ws.Cells.Item(10, 5).Value = '  +1.82%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.7904'
$ws.Cells.Item(11, 5).Value = '  -0.65%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''0.07760'
$ws.Cells.Item(12, 5).Value = '  -0.92%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.822.62'
$ws.Cells.Item(13, 5).Value = '  +0.78%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''90.02'
$ws.Cells.Item(14, 5).Value = '  +4.33%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''5.113'
$ws.Cells.Item(15, 5).Value = '  +2.09%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  -0.05%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '''14.39'
$ws.Cells.Item(17, 5).Value = '  +3.33%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''0.000008018'
$ws.Cells.Item(18, 5).Value = '  +1.16%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -0.15%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '26.801.70'
$ws.Cells.Item(20, 5).Value = '  +4.70%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '2.071.17'
$ws.Cells.Item(21, 5).Value = '  +1.70%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''4.652'
$ws.Cells.Item(22, 5).Value = '  -0.97%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''9.906'
$ws.Cells.Item(23, 5).Value = '  +0.63%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''6.069'
$ws.Cells.Item(24, 5).Value = '  -0.09%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''2.392'
$ws.Cells.Item(25, 5).Value = '  +6.86%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''146.03'
$ws.Cells.Item(26, 5).Value = '  +2.58%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''1.654'
$ws.Cells.Item(27, 5).Value = '  -0.23%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''17.22'
$ws.Cells.Item(28, 5).Value = '  +1.45%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''112.86'
$ws.Cells.Item(29, 5).Value = '  +4.43%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +1.92%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''4.292'
$ws.Cells.Item(31, 5).Value = '  +2.58%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''0.08859'
$ws.Cells.Item(32, 5).Value = '  +1.92%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''0.04956'
$ws.Cells.Item(33, 5).Value = '  +4.01%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''1.158'
$ws.Cells.Item(34, 5).Value = '  +3.18%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''0.7272'
$ws.Cells.Item(35, 5).Value = '  +2.56%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''2.875'
$ws.Cells.Item(36, 5).Value = '  +0.84%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''3.195'
$ws.Cells.Item(37, 5).Value = '  +2.73%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'RenderToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(38, 4).Value = '''2.312'
$ws.Cells.Item(38, 5).Value = '  +0.36%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(39, 4).Value = '''0.01840'
$ws.Cells.Item(39, 5).Value = '  +0.93%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''0.5077'
$ws.Cells.Item(40, 5).Value = '  +1.55%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''116.01'
$ws.Cells.Item(41, 5).Value = '  +0.63%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''0.9245'
$ws.Cells.Item(42, 5).Value = '  -0.52%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''6.146'
$ws.Cells.Item(43, 5).Value = '  +0.21%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''7.991'
$ws.Cells.Item(44, 5).Value = '  +3.30%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''1.000'
$ws.Cells.Item(45, 5).Value = '  -0.10%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''0.4400'
$ws.Cells.Item(46, 5).Value = '  +1.35%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''0.1324'
$ws.Cells.Item(47, 5).Value = '  -1.27%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''9.328'
$ws.Cells.Item(48, 5).Value = '  +2.20%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +0.20%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''0.05930'
$ws.Cells.Item(50, 5).Value = '  +2.10%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''1.460'
$ws.Cells.Item(51, 5).Value = '  -0.24%  '

Write-Host "cryptos list updated"
